$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price rows for "Chirimoya" / "Cultivar IV Región" at
# "Vega Central Mapocho de Santiago" dated 44509 (2021-11-09),
# appended after the last existing data row (116).

$rows = @(
    @{ L="Cuarta";                   M=250; N=1200;  O=1200;  P=1200;  Q="$/kilo (en caja de 15 kilos)"; S=1200; T=1 },
    @{ L="Especial";                 M=330; N=20000; O=20000; P=20000; Q="$/bandeja 8 kilos";             S=2500; T=8 },
    @{ L="Extra (doble especial)";   M=300; N=24000; O=24000; P=24000; Q="$/bandeja 8 kilos";             S=3000; T=8 },
    @{ L="Primera";                  M=280; N=16000; O=16000; P=16000; Q="$/bandeja 8 kilos";             S=2000; T=8 },
    @{ L="Segunda";                  M=350; N=14400; O=14400; P=14400; Q="$/bandeja 8 kilos";             S=1800; T=8 },
    @{ L="Tercera";                  M=220; N=1400;  O=1400;  P=1400;  Q="$/kilo (en caja de 15 kilos)"; S=1400; T=1 }
)

$startRow = 117

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = 9
    $ws.Cells.Item($r, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = 44509
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r - 1, 4).NumberFormat
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100107
    $ws.Cells.Item($r, 8).Value = "Otros"
    $ws.Cells.Item($r, 9).Value = 100107002
    $ws.Cells.Item($r, 10).Value = "Chirimoya"
    $ws.Cells.Item($r, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
}

Write-Host "Added rows 117-122"
